# Sacramento.xlsx roster fix: three player rows had gotten their data
# (No., Pos, Ht, Wt, Birth Date, country flag, Exp, College, bbref url)
# swapped with the adjacent row. Swap each pair back (row index in
# column A is left untouched, only columns B:K move).
#
# The swap is done with Copy/Paste (via a scratch row far below the
# table) rather than reading/writing `.Value`, because assigning a
# numeric-looking string like "5" through `.Value` gets auto-coerced to
# a number (losing the original shared-string/text cell type used by
# the "Exp" column, which mixes numbers with "R" for rookies). Copying
# whole ranges preserves each cell's original type and style exactly,
# and `.Delete()` (rather than `.ClearContents()`) on the scratch row
# removes it completely so no stray row/dimension growth is left
# behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratchRow = 100

function Swap-RosterRows($rowA, $rowB) {
    $rangeA = $ws.Range("B$rowA`:K$rowA")
    $rangeB = $ws.Range("B$rowB`:K$rowB")
    $scratch = $ws.Range("B$scratchRow`:K$scratchRow")

    $rangeA.Copy($scratch)
    $rangeB.Copy($rangeA)
    $scratch.Copy($rangeB)
    $scratch.Delete()
}

# Kevin Huerter (row 6) <-> Malik Monk (row 7)
Swap-RosterRows 6 7

# Trey Lyles (row 8) <-> De'Aaron Fox (row 9)
Swap-RosterRows 8 9

# PJ Dozier (row 16) <-> Kessler Edwards (row 17)
Swap-RosterRows 16 17
